# Daily attendance processing - 2026-01-17 06:44:00
# Normalize the "Recorded By" (column G) values: change "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System" wherever it occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
